$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notifications")

# Row 10 ("On Win ") -> mark as Done and add the Api link, matching the
# other completed rows (e.g. rows 2, 3, 5-8): Status = "Done", Api = "Api".
$ws.Range("E10").Value = "Done"

$ws.Range("F10").Value = "Api"
# Column F has no sheet-level default style, so line up F10's alignment
# with the rest of the "Api" column (centered, like F2/F3/F5:F8).
$ws.Range("F10").VerticalAlignment = -4108
$ws.Range("F10").HorizontalAlignment = -4108

# Move the active selection from F10 to E10.
$ws.Range("E10").Select()
